$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "<the>"
$ws.Range("C2").Value = 30

# Row 3
$ws.Range("C3").Value = 34

# Row 4
$ws.Range("C4").Value = 33

# Row 5
$ws.Range("C5").Value = 31

# Row 6
$ws.Range("B6").Value = "<it>"
$ws.Range("C6").Value = 30

# Row 7
$ws.Range("C7").Value = 28

# Row 8
$ws.Range("B8").Value = "<dumber>"
$ws.Range("C8").Value = 30

# Row 9
$ws.Range("B9").Value = "<foxtrot>"
$ws.Range("C9").Value = 30

# Row 10
$ws.Range("B10").Value = "<a>"
$ws.Range("C10").Value = 25

# Row 11
$ws.Range("B11").Value = "<find>"
$ws.Range("C11").Value = 29

# Row 12
$ws.Range("C12").Value = 31

# Row 13
$ws.Range("C13").Value = 34

# Row 14
$ws.Range("C14").Value = 28

# Row 15
$ws.Range("C15").Value = 31

# Row 16
$ws.Range("B16").Value = "<number>"
$ws.Range("C16").Value = 29

# Row 17
$ws.Range("B17").Value = "<ence>"
$ws.Range("C17").Value = 32

# Row 18
$ws.Range("B18").Value = "<whe>"
$ws.Range("C18").Value = 23
